$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray style from I7 but keep its text ("Out of Sample:") intact.
$ws.Range("I7").ClearFormats()

# Clear the empty "filler" cells that only carried left-over formatting
# (no content at all) from columns H/I/J/Q/S/T. A full Clear() drops the
# cell entirely (contents + formatting) so rows 5, 6 and 27 disappear once
# fully emptied, and the other rows lose their stray styled-but-blank
# cells, matching the cleanup in the diff. (Cleared one address at a time
# -- a single multi-area Range only clears its first area.)
$fillerAddresses = @(
  "I5","I6","J6","J7","Q7","S7","T7","I8","J8","Q8",
  "S8","T8","I9","J9","Q9","S9","T9","I10","J10","Q10",
  "S10","T10","I11","J11","Q11","S11","T11","I12","J12","Q12",
  "S12","T12","I13","J13","Q13","S13","T13","I14","J14","Q14",
  "S14","T14","I15","J15","Q15","S15","T15","I16","J16","Q16",
  "S16","T16","I17","J17","Q17","S17","T17","I18","J18","Q18",
  "S18","T18","I19","J19","Q19","S19","T19","I20","J20","Q20",
  "S20","T20","I21","J21","Q21","S21","T21","I22","J22","Q22",
  "S22","T22","I23","J23","Q23","S23","T23","I24","J24","Q24",
  "S24","T24","I25","J25","Q25","S25","T25","H26","I26","J26",
  "Q26","S26","T26","H27","I27","J27","Q27","H28","I28","J28",
  "Q28","H29","Q29","Q30","Q31"
)
foreach ($addr in $fillerAddresses) {
    $ws.Range($addr).Clear()
}

# Mark Daniel's rows in column R (rows 19-24), and add the final
# analysis/graphs note on row 26.
$ws.Range("R19").Value = "Daniel"
$ws.Range("R20").Value = "Daniel"
$ws.Range("R21").Value = "Daniel"
$ws.Range("R22").Value = "Daniel"
$ws.Range("R23").Value = "Daniel"
$ws.Range("R24").Value = "Daniel"
$ws.Range("R26").Value = "Daniel"
$ws.Range("S26").Value = "Final analysis and graphs"

# Update the view: drop the old frozen/scrolled top-left cell and move the
# selection to R25.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("R25").Select()
